$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 17284
$ws.Range("E2").Value = 1137
$ws.Range("F2").Value = 1137
$ws.Range("G2").Value = 1045
$ws.Range("H2").Value = 778
$ws.Range("I2").Value = 349
$ws.Range("J2").Value = 429
$ws.Range("K2").Value = 11213
$ws.Range("L2").Value = 6751
$ws.Range("M2").Value = 4463
$ws.Range("N2").Value = 2154
$ws.Range("O2").Value = 2308
$ws.Range("P2").Value = 200
$ws.Range("Q2").Value = 1124
$ws.Range("R2").Value = -1102
$ws.Range("S2").Value = 275
$ws.Range("T2").Value = 318
$ws.Range("U2").Value = 806
$ws.Range("V2").Value = 4222
$ws.Range("W2").Value = 6.58
$ws.Range("X2").Value = 4.5
$ws.Range("Y2").Value = 17.55
$ws.Range("Z2").Value = 7.8
$ws.Range("AA2").Value = 151.27
$ws.Range("AB2").Value = 1215.67
$ws.Range("AC2").Value = 872
$ws.Range("AD2").Value = 14.16
$ws.Range("AE2").Value = 5386
$ws.Range("AF2").Value = 2.29
$ws.Range("AG2").Value = 100
$ws.Range("AH2").Value = 0.81
$ws.Range("AI2").Value = 11.46
$ws.Range("AJ2").Value = 40000000

# Row 3
$ws.Range("D3").Value = 20860
$ws.Range("E3").Value = 1628
$ws.Range("F3").Value = 1628
$ws.Range("G3").Value = 1488
$ws.Range("H3").Value = 1121
$ws.Range("I3").Value = 529
$ws.Range("J3").Value = 592
$ws.Range("K3").Value = 14361
$ws.Range("L3").Value = 8869
$ws.Range("M3").Value = 5492
$ws.Range("N3").Value = 2635
$ws.Range("O3").Value = 2857
$ws.Range("P3").Value = 200
$ws.Range("Q3").Value = 418
$ws.Range("R3").Value = -827
$ws.Range("S3").Value = 1294
$ws.Range("T3").Value = 425
$ws.Range("U3").Value = -7
$ws.Range("V3").Value = 5791
$ws.Range("W3").Value = 7.8
$ws.Range("X3").Value = 5.38
$ws.Range("Y3").Value = 22.09
$ws.Range("Z3").Value = 8.77
$ws.Range("AA3").Value = 161.47
$ws.Range("AB3").Value = 1456.93
$ws.Range("AC3").Value = 1323
$ws.Range("AD3").Value = 18.98
$ws.Range("AE3").Value = 6588
$ws.Range("AF3").Value = 3.81
$ws.Range("AG3").Value = 120
$ws.Range("AH3").Value = 0.48
$ws.Range("AI3").Value = 9.07
$ws.Range("AJ3").Value = 40000000

# Row 4
$ws.Range("D4").Value = 22400
$ws.Range("E4").Value = 861
$ws.Range("F4").Value = 861
$ws.Range("G4").Value = 631
$ws.Range("H4").Value = 421
$ws.Range("I4").Value = 155
$ws.Range("J4").Value = 266
$ws.Range("K4").Value = 17406
$ws.Range("L4").Value = 10606
$ws.Range("M4").Value = 6800
$ws.Range("N4").Value = 2649
$ws.Range("O4").Value = 4151
$ws.Range("P4").Value = 200
$ws.Range("Q4").Value = 843
$ws.Range("R4").Value = -672
$ws.Range("S4").Value = 569
$ws.Range("T4").Value = 545
$ws.Range("U4").Value = 298
$ws.Range("V4").Value = 6968
$ws.Range("W4").Value = 3.84
$ws.Range("X4").Value = 1.88
$ws.Range("Y4").Value = 5.88
$ws.Range("Z4").Value = 2.65
$ws.Range("AA4").Value = 155.97
$ws.Range("AB4").Value = 1468.55
$ws.Range("AC4").Value = 388
$ws.Range("AD4").Value = 26.8
$ws.Range("AE4").Value = 6636
$ws.Range("AF4").Value = 1.57
$ws.Range("AG4").Value = 160
$ws.Range("AH4").Value = 1.54
$ws.Range("AI4").Value = 41.18
$ws.Range("AJ4").Value = 40000000

# Row 5
$ws.Range("D5").Value = 24623
$ws.Range("E5").Value = 784
$ws.Range("F5").Value = 784
$ws.Range("G5").Value = 880
$ws.Range("H5").Value = 604
$ws.Range("I5").Value = 317
$ws.Range("J5").Value = 287
$ws.Range("K5").Value = 16936
$ws.Range("L5").Value = 9924
$ws.Range("M5").Value = 7012
$ws.Range("N5").Value = 2853
$ws.Range("O5").Value = 4159
$ws.Range("P5").Value = 200
$ws.Range("Q5").Value = 958
$ws.Range("R5").Value = -1727
$ws.Range("S5").Value = -753
$ws.Range("T5").Value = 349
$ws.Range("U5").Value = 609
$ws.Range("V5").Value = 6382
$ws.Range("W5").Value = 3.19
$ws.Range("X5").Value = 2.45
$ws.Range("Y5").Value = 11.52
$ws.Range("Z5").Value = 3.52
$ws.Range("AA5").Value = 141.53
$ws.Range("AB5").Value = 1586.73
$ws.Range("AC5").Value = 792
$ws.Range("AD5").Value = 13.57
$ws.Range("AE5").Value = 7242
$ws.Range("AF5").Value = 1.48
$ws.Range("AG5").Value = 220
$ws.Range("AH5").Value = 2.05
$ws.Range("AI5").Value = 27.34
$ws.Range("AJ5").Value = 40000000

# Row 6
$ws.Range("D6").Value = 25215
$ws.Range("E6").Value = 565
$ws.Range("F6").Value = 565
$ws.Range("G6").Value = -366
$ws.Range("H6").Value = -454
$ws.Range("I6").Value = -159
$ws.Range("K6").Value = 16783
$ws.Range("L6").Value = 10606
$ws.Range("M6").Value = 6177
$ws.Range("N6").Value = 2604
$ws.Range("P6").Value = 200
$ws.Range("Q6").Value = 145
$ws.Range("R6").Value = -276
$ws.Range("S6").Value = -41
$ws.Range("T6").Value = 407
$ws.Range("U6").Value = -262
$ws.Range("V6").Value = 6736
$ws.Range("W6").Value = 2.24
$ws.Range("X6").Value = -1.8
$ws.Range("Y6").Value = -5.83
$ws.Range("Z6").Value = -2.69
$ws.Range("AA6").Value = 171.69
$ws.Range("AB6").Value = 1495.74
$ws.Range("AC6").Value = -398
$ws.Range("AD6").Value = -20.13
$ws.Range("AE6").Value = 6633
$ws.Range("AF6").Value = 1.21
$ws.Range("AG6").Value = 220
$ws.Range("AH6").Value = 2.75
$ws.Range("AI6").Value = -54.26
$ws.Range("AJ6").Value = 40000000

# Row 7
$ws.Range("D7").Value = 27948
$ws.Range("E7").Value = 1066
$ws.Range("G7").Value = 668
$ws.Range("H7").Value = 504
$ws.Range("I7").Value = 296
$ws.Range("K7").Value = 17582
$ws.Range("L7").Value = 11037
$ws.Range("M7").Value = 6546
$ws.Range("N7").Value = 2879
$ws.Range("P7").Value = 200
$ws.Range("Q7").Value = 1030
$ws.Range("R7").Value = -2426
$ws.Range("S7").Value = -118
$ws.Range("T7").Value = 476
$ws.Range("U7").Value = 1074
$ws.Range("W7").Value = 3.82
$ws.Range("X7").Value = 1.8
$ws.Range("Y7").Value = 10.78
$ws.Range("Z7").Value = 2.93
$ws.Range("AA7").Value = 168.62
$ws.Range("AC7").Value = 739
$ws.Range("AD7").Value = 8.49
$ws.Range("AE7").Value = 7333
$ws.Range("AF7").Value = 0.86
$ws.Range("AG7").Value = 220
$ws.Range("AH7").Value = 3.51
$ws.Range("AI7").Value = 29.78

# Row 8
$ws.Range("D8").Value = 29812
$ws.Range("E8").Value = 1377
$ws.Range("G8").Value = 1082
$ws.Range("H8").Value = 822
$ws.Range("I8").Value = 400
$ws.Range("K8").Value = 18279
$ws.Range("L8").Value = 11223
$ws.Range("M8").Value = 7056
$ws.Range("N8").Value = 3276
$ws.Range("P8").Value = 200
$ws.Range("Q8").Value = 1406
$ws.Range("R8").Value = -2624
$ws.Range("S8").Value = -110
$ws.Range("T8").Value = 412
$ws.Range("U8").Value = -4
$ws.Range("W8").Value = 4.62
$ws.Range("X8").Value = 2.76
$ws.Range("Y8").Value = 13
$ws.Range("Z8").Value = 4.58
$ws.Range("AA8").Value = 159.06
$ws.Range("AC8").Value = 1000
$ws.Range("AD8").Value = 6.27
$ws.Range("AE8").Value = 8344
$ws.Range("AF8").Value = 0.75
$ws.Range("AG8").Value = 220
$ws.Range("AH8").Value = 3.51
$ws.Range("AI8").Value = 22

# Row 9
$ws.Range("D9").Value = 31864
$ws.Range("E9").Value = 1734
$ws.Range("G9").Value = 1470
$ws.Range("H9").Value = 1116
$ws.Range("I9").Value = 536
$ws.Range("K9").Value = 19098
$ws.Range("L9").Value = 11352
$ws.Range("M9").Value = 7746
$ws.Range("N9").Value = 3782
$ws.Range("P9").Value = 200
$ws.Range("Q9").Value = 952
$ws.Range("R9").Value = -2628
$ws.Range("S9").Value = -110
$ws.Range("T9").Value = 404
$ws.Range("U9").Value = -146
$ws.Range("W9").Value = 5.44
$ws.Range("X9").Value = 3.5
$ws.Range("Y9").Value = 15.19
$ws.Range("Z9").Value = 5.97
$ws.Range("AA9").Value = 146.54
$ws.Range("AC9").Value = 1340
$ws.Range("AD9").Value = 4.68
$ws.Range("AE9").Value = 9633
$ws.Range("AF9").Value = 0.65
$ws.Range("AG9").Value = 220
$ws.Range("AH9").Value = 3.51
$ws.Range("AI9").Value = 16.42

